$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must land as TEXT (shared string) in the
# worksheet XML even when it looks like a plain integer (e.g. "2012").
# Assigning such a string straight to .Value lets Excel's COM layer
# auto-coerce it to a numeric cell, so instead we stash it behind a
# formula that evaluates to text and then flatten the formula down to
# its cached value with a values-only paste - this keeps the result a
# plain <c t="s"> shared-string cell with no left-over cell styling.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

# Row 23: institution / dfi - start_date - 2012
$ws.Range("A23").Value = 1747340787
$ws.Range("B23").Value = "update"
$ws.Range("C23").Value = "institution"
$ws.Range("D23").Value = "dfi"
$ws.Range("F23").Value = "start_date"
Set-TextValue $ws.Range("H23") "2012"

# Row 24: institution / dfi - end_date - 2023
$ws.Range("A24").Value = 1747340787
$ws.Range("B24").Value = "update"
$ws.Range("C24").Value = "institution"
$ws.Range("D24").Value = "dfi"
$ws.Range("F24").Value = "end_date"
Set-TextValue $ws.Range("H24") "2023"

# Row 25: institution / dff - start_date - 2010/10
$ws.Range("A25").Value = 1747340787
$ws.Range("B25").Value = "update"
$ws.Range("C25").Value = "institution"
$ws.Range("D25").Value = "dff"
$ws.Range("F25").Value = "start_date"
$ws.Range("H25").Value = "2010/10"

# Row 26: institution / seco - end_date - 2021/04
$ws.Range("A26").Value = 1747340787
$ws.Range("B26").Value = "update"
$ws.Range("C26").Value = "institution"
$ws.Range("D26").Value = "seco"
$ws.Range("F26").Value = "end_date"
$ws.Range("H26").Value = "2021/04"
